$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# "Polarity" column (B) to C and "Review" column (C) to D.
$ws.Columns("B").Insert()

# New column header (mirrors column A, "Unnamed: 0")
$ws.Range("B1").Value = "Unnamed: 0.1"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").Borders.LineStyle = 1

# Fill the new column with the same values as column A (the old index)
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4

# Lowercase the review text, now living in column D
$ws.Range("D2").Value = $ws.Range("D2").Value().ToLower()
$ws.Range("D3").Value = $ws.Range("D3").Value().ToLower()
$ws.Range("D4").Value = $ws.Range("D4").Value().ToLower()
$ws.Range("D5").Value = $ws.Range("D5").Value().ToLower()
$ws.Range("D6").Value = $ws.Range("D6").Value().ToLower()
